# Auto-generated Excel COM-interop script to apply the Halicarnassus_Profits.xlsx diff
# The diff updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ figures across several of the
# per-job Leve profit sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR). Some rows also gain or
# lose a profit cell entirely (where the source data produced/no-longer-produced a value).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 118
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H28").Value = 2832.476
$ws.Range("I28").Value = 733.0769
$ws.Range("K28").Value = 733.0769
$ws.Range("M28").Value = -248.0769
$ws.Range("H80").Value = 600
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 600
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H86").Value = 5748.2144
$ws.Range("I86").Value = 4590.273
$ws.Range("J86").Value = 9994
$ws.Range("K86").Value = 4590.273
$ws.Range("L86").Value = 9994
$ws.Range("M86").Value = -3467.273
$ws.Range("N86").Value = -12240
$ws.Range("H89").Value = 5748.2144
$ws.Range("I89").Value = 4590.273
$ws.Range("J89").Value = 9994
$ws.Range("K89").Value = 22951.365
$ws.Range("L89").Value = 49970
$ws.Range("M89").Value = -17335.365
$ws.Range("N89").Value = -61202
$ws.Range("H92").Value = 164.75
$ws.Range("I92").Value = 164.75
$ws.Range("K92").Value = 164.75
$ws.Range("M92").Value = 1083.25
$ws.Range("H107").Value = 1561.909
$ws.Range("I107").Value = 1880.2222
$ws.Range("J107").Value = 129.5
$ws.Range("K107").Value = 1880.2222
$ws.Range("L107").Value = 129.5
$ws.Range("M107").Value = 39.77780000000007
$ws.Range("N107").Value = -3969.5
$ws.Range("H113").Value = 10000
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H129").Value = 1263.2858
$ws.Range("I129").Value = 973.8333
$ws.Range("K129").Value = 2921.4999
$ws.Range("M129").Value = 2078.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 193.81818
$ws.Range("I4").Value = 102.52631
$ws.Range("K4").Value = 102.52631
$ws.Range("M4").Value = 13.47369
$ws.Range("H6").Value = 16363863
$ws.Range("I6").Value = 20000082
$ws.Range("K6").Value = 20000082
$ws.Range("M6").Value = -19999909
$ws.Range("H12").Value = 10000
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H23").Value = 12749.75
$ws.Range("J23").Value = 12749.75
$ws.Range("L23").Value = 12749.75
$ws.Range("N23").Value = -13267.75
$ws.Range("H26").Value = 707
$ws.Range("I26").Value = 707
$ws.Range("K26").Value = 707
$ws.Range("M26").Value = -377
$ws.Range("H38").Value = 17634.5
$ws.Range("I38").Value = 2505.6667
$ws.Range("J38").Value = 63021
$ws.Range("K38").Value = 2505.6667
$ws.Range("L38").Value = 63021
$ws.Range("M38").Value = -2038.6667
$ws.Range("N38").Value = -63955
$ws.Range("H39").Value = 4638.3335
$ws.Range("I39").Value = 4638.3335
$ws.Range("K39").Value = 4638.3335
$ws.Range("M39").Value = -4118.3335
$ws.Range("H43").Value = 11999994
$ws.Range("J43").Value = 8999987
$ws.Range("L43").Value = 8999987
$ws.Range("N43").Value = -9000613
$ws.Range("H80").Value = 47736.8
$ws.Range("J80").Value = 47736.8
$ws.Range("L80").Value = 47736.8
$ws.Range("N80").Value = -49732.8
$ws.Range("H83").Value = 47736.8
$ws.Range("J83").Value = 47736.8
$ws.Range("L83").Value = 143210.4
$ws.Range("N83").Value = -153194.4
$ws.Range("H114").Value = 7525000
$ws.Range("J114").Value = 7525000
$ws.Range("L114").Value = 7525000
$ws.Range("N114").Value = -7533678
$ws.Range("H132").Value = 5704.091
$ws.Range("I132").Value = 5704.091
$ws.Range("K132").Value = 17112.273
$ws.Range("M132").Value = -14582.273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1804.3334
$ws.Range("I22").Value = 1804.3334
$ws.Range("K22").Value = 1804.3334
$ws.Range("M22").Value = -1631.3334
$ws.Range("H86").Value = 3569.9583
$ws.Range("I86").Value = 1741.4286
$ws.Range("J86").Value = 6129.9
$ws.Range("K86").Value = 1741.4286
$ws.Range("L86").Value = 6129.9
$ws.Range("M86").Value = -618.4286
$ws.Range("N86").Value = -8375.9
$ws.Range("H89").Value = 3569.9583
$ws.Range("I89").Value = 1741.4286
$ws.Range("J89").Value = 6129.9
$ws.Range("K89").Value = 8707.143
$ws.Range("L89").Value = 30649.5
$ws.Range("M89").Value = -3091.143
$ws.Range("N89").Value = -41881.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 51334.332
$ws.Range("I6").Value = 40000
$ws.Range("K6").Value = 40000
$ws.Range("M6").Value = -39887
$ws.Range("H7").Value = 3453
$ws.Range("I7").Value = 3928.6538
$ws.Range("K7").Value = 3928.6538
$ws.Range("M7").Value = -3815.6538
$ws.Range("H31").Value = 4399.976
$ws.Range("I31").Value = 2063.1035
$ws.Range("K31").Value = 2063.1035
$ws.Range("M31").Value = -1768.1035
$ws.Range("H32").Value = 1801.5
$ws.Range("I32").Value = 1801.5
$ws.Range("K32").Value = 1801.5
$ws.Range("M32").Value = -1485.5
$ws.Range("H34").Value = 4399.976
$ws.Range("I34").Value = 2063.1035
$ws.Range("K34").Value = 2063.1035
$ws.Range("M34").Value = -1861.1035
$ws.Range("H99").Value = 2400
$ws.Range("J99").Value = 2400
$ws.Range("L99").Value = 2400
$ws.Range("N99").Value = -5396
$ws.Range("H126").Value = 2400
$ws.Range("J126").Value = 2400
$ws.Range("L126").Value = 7200
$ws.Range("N126").Value = -12140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 278.91666
$ws.Range("I2").Value = 176.6
$ws.Range("K2").Value = 176.6
$ws.Range("M2").Value = -63.59999999999999
$ws.Range("H80").Value = 3430
$ws.Range("I80").Value = 3395
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 3395
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = -2397
$ws.Range("N80").Value = -5496
$ws.Range("H83").Value = 3430
$ws.Range("I83").Value = 3395
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 16975
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = -11983
$ws.Range("N83").Value = -27484
$ws.Range("H132").Value = 2321.7778
$ws.Range("I132").Value = 1985.1428
$ws.Range("K132").Value = 5955.428400000001
$ws.Range("M132").Value = -3425.428400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 3000
$ws.Range("I19").Value = 500
$ws.Range("J19").Value = 3833.3333
$ws.Range("K19").Value = 500
$ws.Range("L19").Value = 3833.3333
$ws.Range("M19").Value = -330
$ws.Range("N19").Value = -4173.3333
$ws.Range("H32").Value = 17506.5
$ws.Range("I32").Value = 17506.5
$ws.Range("K32").Value = 17506.5
$ws.Range("M32").Value = -17189.5
$ws.Range("H46").Value = 6497.5
$ws.Range("I46").Value = 1997.5
$ws.Range("K46").Value = 1997.5
$ws.Range("M46").Value = -1809.5
$ws.Range("H100").Value = 7200
$ws.Range("I100").Value = 3500
$ws.Range("J100").Value = 9666.666999999999
$ws.Range("K100").Value = 3500
$ws.Range("L100").Value = 9666.666999999999
$ws.Range("M100").Value = -2959
$ws.Range("N100").Value = -10748.667
$ws.Range("H116").Value = 184000
$ws.Range("J116").Value = 184000
$ws.Range("L116").Value = 184000
$ws.Range("N116").Value = -193178
$ws.Range("H122").Value = 2687
$ws.Range("I122").Value = 2326.7144
$ws.Range("K122").Value = 6980.1432
$ws.Range("M122").Value = -4530.1432
$ws.Range("H132").Value = 3923.75
$ws.Range("I132").Value = 2977
$ws.Range("J132").Value = 5501.6665
$ws.Range("K132").Value = 8931
$ws.Range("L132").Value = 16504.9995
$ws.Range("M132").Value = -6401
$ws.Range("N132").Value = -21564.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2291.111
$ws.Range("I132").Value = 2088.7144
$ws.Range("K132").Value = 6266.1432
$ws.Range("M132").Value = -3736.1432
